$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.170519828796387
$ws.Range("B1").Value = 7.307291507720947
$ws.Range("C1").Value = 5.626019477844238
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 3.754376649856567
